$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, thin border) from an existing
# header cell onto the newly-added header columns (L1:W1) before writing
# their text, so the new header cells pick up style index 1 like B1:K1.
$ws.Range("A1").Copy()
$ws.Range("L1:W1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row (row 1)
$ws.Range("A1").Value = 'productIds'
$ws.Range("B1").Value = 'MSE_no_transfer'
$ws.Range("C1").Value = 'MSE_transfer_basic'
$ws.Range("D1").Value = 'MSE_transfer_coral'
$ws.Range("E1").Value = 'MSE_transfer_sa'
$ws.Range("F1").Value = 'MSE_transfer_bw'
$ws.Range("G1").Value = 'MSE_transfer_nnw'
$ws.Range("H1").Value = 'MAE_no_transfer'
$ws.Range("I1").Value = 'MAE_transfer_basic'
$ws.Range("J1").Value = 'MAE_transfer_coral'
$ws.Range("K1").Value = 'MAE_transfer_sa'
$ws.Range("L1").Value = 'MAE_transfer_bw'
$ws.Range("M1").Value = 'MAE_transfer_nnw'
$ws.Range("N1").Value = 'MSE_diff_basic'
$ws.Range("O1").Value = 'MSE_transfer_coral'
$ws.Range("P1").Value = 'MSE_diff_sa'
$ws.Range("Q1").Value = 'MSE_diff_bw'
$ws.Range("R1").Value = 'MSE_diff_nnw'
$ws.Range("S1").Value = 'MAE_diff_basic'
$ws.Range("T1").Value = 'MAE_transfer_coral'
$ws.Range("U1").Value = 'MAE_diff_sa'
$ws.Range("V1").Value = 'MAE_diff_bw'
$ws.Range("W1").Value = 'MAE_diff_nnw'

# Data rows 2-6
$ws.Range("A2").Value = '101-120'
$ws.Range("B2").Value = 1.222218837026467
$ws.Range("C2").Value = 1.063419466375333
$ws.Range("D2").Value = 1.111192752731464
$ws.Range("E2").Value = 1.020436886935692
$ws.Range("F2").Value = 0.8775133364974671
$ws.Range("G2").Value = 0.9284918576178334
$ws.Range("H2").Value = 0.5127839979635104
$ws.Range("I2").Value = 0.6062345606259164
$ws.Range("J2").Value = 0.5632716150102758
$ws.Range("K2").Value = 0.7701338342334212
$ws.Range("L2").Value = 0.4411453296771398
$ws.Range("M2").Value = 0.4037414647687326
$ws.Range("N2").Value = -0.1587993706511339
$ws.Range("O2").Value = -0.111026084295003
$ws.Range("P2").Value = -0.2017819500907749
$ws.Range("Q2").Value = -0.3447055005290002
$ws.Range("R2").Value = -0.2937269794086339
$ws.Range("S2").Value = 0.09345056266240603
$ws.Range("T2").Value = 0.05048761704676541
$ws.Range("U2").Value = 0.2573498362699108
$ws.Range("V2").Value = -0.07163866828637061
$ws.Range("W2").Value = -0.1090425331947777
$ws.Range("A3").Value = '121-140'
$ws.Range("B3").Value = 2.468190564009104
$ws.Range("C3").Value = 2.589472281019797
$ws.Range("D3").Value = 3.066289370625425
$ws.Range("E3").Value = 2.828577720844734
$ws.Range("F3").Value = 2.769117160019671
$ws.Range("G3").Value = 2.712442318767745
$ws.Range("H3").Value = 0.9699308876351251
$ws.Range("I3").Value = 0.9062258530067143
$ws.Range("J3").Value = 0.8561641027395959
$ws.Range("K3").Value = 1.018480999614012
$ws.Range("L3").Value = 0.8530154002302385
$ws.Range("M3").Value = 0.751503267021855
$ws.Range("N3").Value = 0.1212817170106932
$ws.Range("O3").Value = 0.5980988066163211
$ws.Range("P3").Value = 0.3603871568356301
$ws.Range("Q3").Value = 0.3009265960105671
$ws.Range("R3").Value = 0.2442517547586416
$ws.Range("S3").Value = -0.06370503462841082
$ws.Range("T3").Value = -0.1137667848955293
$ws.Range("U3").Value = 0.04855011197888692
$ws.Range("V3").Value = -0.1169154874048867
$ws.Range("W3").Value = -0.2184276206132701
$ws.Range("A4").Value = '141-160'
$ws.Range("B4").Value = 5.713082570197582
$ws.Range("C4").Value = 5.214058498148367
$ws.Range("D4").Value = 3.489666168914155
$ws.Range("E4").Value = 3.13657315087518
$ws.Range("F4").Value = 3.189241431831304
$ws.Range("G4").Value = 3.116256800558963
$ws.Range("H4").Value = 1.292643636471448
$ws.Range("I4").Value = 1.285187414189988
$ws.Range("J4").Value = 1.004404575526774
$ws.Range("K4").Value = 1.09710563229861
$ws.Range("L4").Value = 0.8976243119413795
$ws.Range("M4").Value = 0.8794305171610247
$ws.Range("N4").Value = -0.499024072049215
$ws.Range("O4").Value = -2.223416401283427
$ws.Range("P4").Value = -2.576509419322401
$ws.Range("Q4").Value = -2.523841138366278
$ws.Range("R4").Value = -2.596825769638619
$ws.Range("S4").Value = -0.007456222281459146
$ws.Range("T4").Value = -0.2882390609446737
$ws.Range("U4").Value = -0.1955380041728372
$ws.Range("V4").Value = -0.3950193245300681
$ws.Range("W4").Value = -0.4132131193104228
$ws.Range("A5").Value = '161-180'
$ws.Range("B5").Value = 13.78807970264317
$ws.Range("C5").Value = 13.45158414717129
$ws.Range("D5").Value = 12.51545602290823
$ws.Range("E5").Value = 12.71251798689346
$ws.Range("F5").Value = 12.58048641638751
$ws.Range("G5").Value = 12.53523161432398
$ws.Range("H5").Value = 1.097150881057269
$ws.Range("I5").Value = 1.164234795733726
$ws.Range("J5").Value = 1.013350874587507
$ws.Range("K5").Value = 1.018628110968516
$ws.Range("L5").Value = 1.000387708345127
$ws.Range("M5").Value = 0.9771102808825743
$ws.Range("N5").Value = -0.3364955554718883
$ws.Range("O5").Value = -1.272623679734945
$ws.Range("P5").Value = -1.075561715749712
$ws.Range("Q5").Value = -1.207593286255666
$ws.Range("R5").Value = -1.252848088319194
$ws.Range("S5").Value = 0.06708391467645725
$ws.Range("T5").Value = -0.08380000646976216
$ws.Range("U5").Value = -0.07852277008875319
$ws.Range("V5").Value = -0.09676317271214208
$ws.Range("W5").Value = -0.1200406001746944
$ws.Range("A6").Value = '181-200'
$ws.Range("B6").Value = 8.304673529463534
$ws.Range("C6").Value = 7.668290373907366
$ws.Range("D6").Value = 5.896055479847965
$ws.Range("E6").Value = 5.820702221270763
$ws.Range("F6").Value = 5.856413974995616
$ws.Range("G6").Value = 6.142760710462269
$ws.Range("H6").Value = 1.053675744853724
$ws.Range("I6").Value = 1.063607334170416
$ws.Range("J6").Value = 0.8640754600180854
$ws.Range("K6").Value = 0.822392764717888
$ws.Range("L6").Value = 0.808529699177241
$ws.Range("M6").Value = 0.7930880846202476
$ws.Range("N6").Value = -0.6363831555561674
$ws.Range("O6").Value = -2.408618049615568
$ws.Range("P6").Value = -2.48397130819277
$ws.Range("Q6").Value = -2.448259554467918
$ws.Range("R6").Value = -2.161912819001264
$ws.Range("S6").Value = 0.009931589316691358
$ws.Range("T6").Value = -0.1896002848356387
$ws.Range("U6").Value = -0.2312829801358361
$ws.Range("V6").Value = -0.2451460456764831
$ws.Range("W6").Value = -0.2605876602334766
